$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 48; this shifts rows 48-67
# down to 49-68 (dates, volumes, prices, origin, etc. all move with
# their row, matching the diff exactly).
$ws.Rows.Item(48).Insert()

# Fill in the brand-new row 48 with its own data (a new weekly entry).
$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = 44988
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = "Poroto granado"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 80
$ws.Range("K48").Value = 36000
$ws.Range("L48").Value = 36000
$ws.Range("M48").Value = 36000
$ws.Range("N48").Value = "`$/saco 25 kilos"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 1440
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
